$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): B3 value corrected, D3 cleared (no longer has a value)
$ws.Range("B3").Value = 237176.5490839333
$ws.Range("D3").ClearContents()

# Row 4 (Methanol): C4 value corrected
$ws.Range("C4").Value = 184.7257177872526

# Row 5 (Ammonia): C5 value corrected
$ws.Range("C5").Value = 5544.926562293677

# Row 7: label changed from "Other" to "Biogas", D7 value corrected
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 1843.392663657924

# New row 8: "Other" row with corrected D8 value
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 1389.575972385624
